$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price (D) cells whose new values look numeric,
# so Excel keeps them as literal text (matching the source data) instead
# of auto-converting to a Number and losing formatting like trailing zeros.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price / Volume(1h) values scraped for this run.
$ws.Range("D2").Value = '27.543.32'
$ws.Range("E2").Value = '  +6.01%  '
$ws.Range("D3").Value = '1.814.06'
$ws.Range("E3").Value = '  +5.81%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '344.15'
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("E7").Value = '  +4.22%  '
$ws.Range("D8").Value = '0.3518'
$ws.Range("E8").Value = '  +5.59%  '
$ws.Range("D9").Value = '49.72'
$ws.Range("E9").Value = '  +1.48%  '
$ws.Range("D10").Value = '1.233'
$ws.Range("E10").Value = '  +4.53%  '
$ws.Range("D11").Value = '0.07762'
$ws.Range("E11").Value = '  +3.55%  '
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").Value = '22.40'
$ws.Range("E13").Value = '  +11.32%  '
$ws.Range("D14").Value = '6.611'
$ws.Range("E14").Value = '  +6.00%  '
$ws.Range("D15").Value = '7.216'
$ws.Range("E15").Value = '  +4.59%  '
$ws.Range("D16").Value = '1.812.92'
$ws.Range("E16").Value = '  +5.89%  '
$ws.Range("D17").Value = '0.00001128'
$ws.Range("E17").Value = '  +5.03%  '
$ws.Range("D18").Value = '0.06724'
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").Value = '86.47'
$ws.Range("E19").Value = '  +5.61%  '
$ws.Range("D20").Value = '1.0000'
$ws.Range("D21").Value = '17.77'
$ws.Range("E21").Value = '  +9.00%  '
$ws.Range("D22").Value = '6.524'
$ws.Range("E22").Value = '  +7.33%  '
$ws.Range("D23").Value = '13.14'
$ws.Range("E23").Value = '  +0.72%  '
$ws.Range("D24").Value = '27.525.73'
$ws.Range("E24").Value = '  +6.04%  '
$ws.Range("D25").Value = '2.469'
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").Value = '2.678'
$ws.Range("E26").Value = '  +7.27%  '
$ws.Range("D27").Value = '22.12'
$ws.Range("E27").Value = '  +14.76%  '
$ws.Range("D28").Value = '1.491'
$ws.Range("E28").Value = '  +14.66%  '
$ws.Range("D29").Value = '153.79'
$ws.Range("E29").Value = '  +2.56%  '
$ws.Range("D30").Value = '2.016.38'
$ws.Range("E30").Value = '  +6.25%  '
$ws.Range("D31").Value = '136.44'
$ws.Range("E31").Value = '  +5.65%  '
$ws.Range("D32").Value = '6.378'
$ws.Range("E32").Value = '  +6.74%  '
$ws.Range("D33").Value = '4.084'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '14.01'
$ws.Range("E34").Value = '  +8.53%  '
$ws.Range("D35").Value = '0.08825'
$ws.Range("E35").Value = '  +3.66%  '
$ws.Range("D36").Value = '1.719'
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("D37").Value = '5.647'
$ws.Range("E37").Value = '  +5.14%  '
$ws.Range("D38").Value = '0.7096'
$ws.Range("E38").Value = '  +15.13%  '
$ws.Range("E39").Value = '  +6.39%  '
$ws.Range("D40").Value = '0.06550'
$ws.Range("E40").Value = '  +5.38%  '
$ws.Range("D41").Value = '0.2273'
$ws.Range("D42").Value = '9.026'
$ws.Range("E42").Value = '  +5.81%  '
$ws.Range("D43").Value = '1.288'
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = '14.88'
$ws.Range("E44").Value = '  +1.56%  '
$ws.Range("D45").Value = '0.6636'
$ws.Range("E45").Value = '  +12.88%  '
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("D47").Value = '3.962'
$ws.Range("E47").Value = '  +3.15%  '
$ws.Range("D48").Value = '2.184'
$ws.Range("E48").Value = '  +8.60%  '
$ws.Range("D49").Value = '133.22'
$ws.Range("E49").Value = '  +4.99%  '
$ws.Range("D50").Value = '0.07362'
$ws.Range("E50").Value = '  +1.37%  '
$ws.Range("D51").Value = '80.75'
$ws.Range("E51").Value = '  +4.86%  '
